$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("class_entities")
$ws.Activate()
